$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.15967521690695
$ws.Range("C2").Value = 6.498337351986359
$ws.Range("E2").Value = 10.50155289186741
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.709888699607193
$ws.Range("K2").Value = 13.61617721057028
$ws.Range("M2").Value = 15.94939587996503
$ws.Range("N2").Value = 22.88384567116072

$ws.Range("B3").Value = 13.94201437529212
$ws.Range("C3").Value = 6.328519332001004
$ws.Range("E3").Value = 10.27982848283818
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.713070515920299
$ws.Range("K3").Value = 13.4754110540914
$ws.Range("M3").Value = 15.80660041038779
$ws.Range("N3").Value = 22.90305095219939

$ws.Range("B4").Value = 13.81093556697714
$ws.Range("C4").Value = 6.22412746346804
$ws.Range("E4").Value = 10.14483235605414
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.715123786374591
$ws.Range("K4").Value = 13.39207869316383
$ws.Range("M4").Value = 15.72277006638242
$ws.Range("N4").Value = 22.91655504757985

$ws.Range("B5").Value = 13.75823630635322
$ws.Range("C5").Value = 6.181627689541162
$ws.Range("E5").Value = 10.09019219513919
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.715985655712127
$ws.Range("K5").Value = 13.35893549620641
$ws.Range("M5").Value = 15.68960882951189
$ws.Range("N5").Value = 22.92248780305093

$ws.Range("B6").Value = 13.74953096064153
$ws.Range("C6").Value = 6.174575239554905
$ws.Range("E6").Value = 10.08114417401878
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.716130289915296
$ws.Range("K6").Value = 13.35348236678546
$ws.Range("M6").Value = 15.68416377144242
$ws.Range("N6").Value = 22.92349886435011

$ws.Range("B7").Value = 13.81022185243385
$ws.Range("C7").Value = 6.223554035067251
$ws.Range("E7").Value = 10.1440938421663
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.715135307911995
$ws.Range("K7").Value = 13.39162836430854
$ws.Range("M7").Value = 15.72231875008258
$ws.Range("N7").Value = 22.91663331995008

$ws.Range("B8").Value = 14.08413628387144
$ws.Range("C8").Value = 6.439851777999389
$ws.Range("E8").Value = 10.42491526240716
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.710965170027301
$ws.Range("K8").Value = 13.56702018831278
$ws.Range("M8").Value = 15.8993820794473
$ws.Range("N8").Value = 22.8901118753756

$ws.Range("B9").Value = 14.63846930566908
$ws.Range("C9").Value = 6.860034793033652
$ws.Range("E9").Value = 10.98103781686464
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.703573701049371
$ws.Range("K9").Value = 13.93387793436265
$ws.Range("M9").Value = 16.27557697709206
$ws.Range("N9").Value = 22.85172196024397

$ws.Range("B10").Value = 15.05196208423076
$ws.Range("C10").Value = 7.162549703468432
$ws.Range("E10").Value = 11.38820845466112
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.698616423061007
$ws.Range("K10").Value = 14.21505023127812
$ws.Range("M10").Value = 16.56744230011838
$ws.Range("N10").Value = 22.83186779692783

$ws.Range("B11").Value = 15.2405347447559
$ws.Range("C11").Value = 7.298109584011277
$ws.Range("E11").Value = 11.57219414604334
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.696462698651155
$ws.Range("K11").Value = 14.34497841163626
$ws.Range("M11").Value = 16.70308744717518
$ws.Range("N11").Value = 22.82465918318768

$ws.Range("B12").Value = 15.31193838845564
$ws.Range("C12").Value = 7.349093141410846
$ws.Range("E12").Value = 11.64161350742968
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.695661616846899
$ws.Range("K12").Value = 14.39442605822699
$ws.Range("M12").Value = 16.75482350799104
$ws.Range("N12").Value = 22.82219242034726

$ws.Range("B13").Value = 15.29656167964724
$ws.Range("C13").Value = 7.338129306770417
$ws.Range("E13").Value = 11.62667514648547
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.695833501259355
$ws.Range("K13").Value = 14.38376635239728
$ws.Range("M13").Value = 16.74366544815967
$ws.Range("N13").Value = 22.82271197281861

$ws.Range("B14").Value = 15.24640966858007
$ws.Range("C14").Value = 7.302311342894084
$ws.Range("E14").Value = 11.57791077796253
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.696396503341591
$ws.Range("K14").Value = 14.34904180765765
$ws.Range("M14").Value = 16.70733663333365
$ws.Range("N14").Value = 22.82445096595595

$ws.Range("B15").Value = 15.2156873787866
$ws.Range("C15").Value = 7.28032468289551
$ws.Range("E15").Value = 11.54800625133329
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.69674324234542
$ws.Range("K15").Value = 14.32780283883329
$ws.Range("M15").Value = 16.68513110785387
$ws.Range("N15").Value = 22.82555042146799

$ws.Range("B16").Value = 15.03964166375193
$ws.Range("C16").Value = 7.153644426833304
$ws.Range("E16").Value = 11.37615321944436
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.698759206108328
$ws.Range("K16").Value = 14.20659603300862
$ws.Range("M16").Value = 16.55863173135874
$ws.Range("N16").Value = 22.83237565338685

$ws.Range("B17").Value = 14.93171197342043
$ws.Range("C17").Value = 7.075363546695083
$ws.Range("E17").Value = 11.27035582579481
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 3.700021833131864
$ws.Range("K17").Value = 14.13272601809513
$ws.Range("M17").Value = 16.48173408526431
$ws.Range("N17").Value = 22.83703025051509

$ws.Range("B18").Value = 14.86968238892244
$ws.Range("C18").Value = 7.030147937063358
$ws.Range("E18").Value = 11.20939157378107
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 3.700757609195337
$ws.Range("K18").Value = 14.09043085529212
$ws.Range("M18").Value = 16.4377777415586
$ws.Range("N18").Value = 22.83987900716037

$ws.Range("B19").Value = 14.84869085629589
$ws.Range("C19").Value = 7.014807718404757
$ws.Range("E19").Value = 11.18873314609371
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 3.701008372660445
$ws.Range("K19").Value = 14.07614491478166
$ws.Range("M19").Value = 16.42294309795017
$ws.Range("N19").Value = 22.84087298677789

$ws.Range("B20").Value = 14.94319676487803
$ws.Range("C20").Value = 7.083716806317613
$ws.Range("E20").Value = 11.28163032019233
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("G20").Value = 3.699886436966943
$ws.Range("K20").Value = 14.14056996029469
$ws.Range("M20").Value = 16.4898920120356
$ws.Range("N20").Value = 22.83651699952553

$ws.Range("B21").Value = 15.26114122634387
$ws.Range("C21").Value = 7.312841858361861
$ws.Range("E21").Value = 11.59224146302405
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.696230743523937
$ws.Range("K21").Value = 14.35923490425736
$ws.Range("M21").Value = 16.71799758571498
$ws.Range("N21").Value = 22.82393303782687

$ws.Range("B22").Value = 15.46887672036499
$ws.Range("C22").Value = 7.460525997714193
$ws.Range("E22").Value = 11.79374358578304
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.693925935944658
$ws.Range("K22").Value = 14.50356252760002
$ws.Range("M22").Value = 16.86921527298117
$ws.Range("N22").Value = 22.81724192087796

$ws.Range("B23").Value = 15.35803323767747
$ws.Range("C23").Value = 7.381909723989836
$ws.Range("E23").Value = 11.68635868804584
$ws.Range("F23").Value = 21.82633154475857
$ws.Range("G23").Value = 3.695148362031828
$ws.Range("K23").Value = 14.42641697403147
$ws.Range("M23").Value = 16.78832630082338
$ws.Range("N23").Value = 22.82067254059109

$ws.Range("B24").Value = 14.93800442082042
$ws.Range("C24").Value = 7.079940950339335
$ws.Range("E24").Value = 11.27653355147487
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.699947618806893
$ws.Range("K24").Value = 14.13702316963522
$ws.Range("M24").Value = 16.48620302158586
$ws.Range("N24").Value = 22.83674850215701

$ws.Range("B25").Value = 14.48710864972838
$ws.Range("C25").Value = 6.747194165981433
$ws.Range("E25").Value = 10.83050893777495
$ws.Range("F25").Value = 18.34778573295691
$ws.Range("G25").Value = 3.705489750226945
$ws.Range("K25").Value = 13.8324164424631
$ws.Range("M25").Value = 16.17092475897439
$ws.Range("N25").Value = 22.86064439224615
